$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "record_atd" (C) and "average_simulation_TD" (D) values
# for Appenzeller-Herzog (2019) - van Dis (2020) relevance markers correction.

$ws.Range("C2").Value = 18
$ws.Range("D2").Value = 16

$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 23

$ws.Range("C4").Value = 29
$ws.Range("D4").Value = 18.5

$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 22

$ws.Range("C6").Value = 199
$ws.Range("D6").Value = 194

$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 3

$ws.Range("C8").Value = 49
$ws.Range("D8").Value = 49.5

$ws.Range("C9").Value = 4
$ws.Range("D9").Value = 3.5

$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 18.5

$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 28

$ws.Range("C12").Value = 21
$ws.Range("D12").Value = 19

$ws.Range("C13").Value = 62
$ws.Range("D13").Value = 62

$ws.Range("C14").Value = 22
$ws.Range("D14").Value = 13

$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 18

$ws.Range("C16").Value = 41
$ws.Range("D16").Value = 45

$ws.Range("C17").Value = 35.8
